$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($ws, $row, $vals)
    foreach ($col in $vals.Keys) {
        $ws.Cells.Item($row, $col).Value2 = $vals[$col]
    }
}

# Swap F:V between row 98 and row 99
$row98 = @{
    6 = 'Al Feiha'
    7 = 1
    8 = 'Al Nassr'
    9 = 3
    10 = 5.59
    11 = '24/10/2023 22:01'
    12 = 6.97
    13 = '28/10/2023 16:59'
    14 = 4.98
    15 = '24/10/2023 22:01'
    16 = 5.6
    17 = '28/10/2023 16:59'
    18 = 1.49
    19 = '24/10/2023 22:01'
    20 = 1.37
    21 = '28/10/2023 16:59'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-feiha-al-nassr/ve12DxUP/'
}
$row99 = @{
    6 = 'Abha'
    7 = 2
    8 = 'Al Shabab'
    9 = 1
    10 = 3.77
    11 = '25/10/2023 13:48'
    12 = 4.27
    13 = '28/10/2023 16:56'
    14 = 3.63
    15 = '25/10/2023 13:48'
    16 = 4.04
    17 = '28/10/2023 16:56'
    18 = 1.87
    19 = '25/10/2023 13:48'
    20 = 1.76
    21 = '28/10/2023 16:56'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/abha-al-shabab/tAraEIEJ/'
}
Set-RowValues $ws 98 $row99
Set-RowValues $ws 99 $row98

# Swap F:V between row 112 and row 113
$row112 = @{
    6 = 'Damac'
    7 = 2
    8 = 'Al Ahli SC'
    9 = 2
    10 = 4.18
    11 = '06/11/2023 03:42'
    12 = 4.06
    13 = '09/11/2023 18:50'
    14 = 4.22
    15 = '06/11/2023 03:42'
    16 = 3.95
    17 = '09/11/2023 18:50'
    18 = 1.74
    19 = '06/11/2023 03:42'
    20 = 1.83
    21 = '09/11/2023 18:50'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/damac-al-ahli-sc/vuLPGBzI/'
}
$row113 = @{
    6 = 'Al Riyadh'
    7 = 1
    8 = 'Al Fateh'
    9 = 1
    10 = 4.45
    11 = '05/11/2023 19:12'
    12 = 4.56
    13 = '09/11/2023 18:58'
    14 = 4.11
    15 = '05/11/2023 19:12'
    16 = 4.22
    17 = '09/11/2023 18:58'
    18 = 1.72
    19 = '05/11/2023 19:12'
    20 = 1.69
    21 = '09/11/2023 18:58'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-riyadh-al-fateh/WCxbb3Kj/'
}
Set-RowValues $ws 112 $row113
Set-RowValues $ws 113 $row112

# Swap F:V between row 116 and row 117
$row116 = @{
    6 = 'Al Akhdoud'
    7 = 2
    8 = 'Al Hazem'
    9 = 1
    10 = 2.19
    11 = '08/11/2023 06:12'
    12 = 1.72
    13 = '11/11/2023 15:56'
    14 = 3.47
    15 = '08/11/2023 06:12'
    16 = 4.09
    17 = '11/11/2023 15:56'
    18 = 3.06
    19 = '08/11/2023 06:12'
    20 = 4.49
    21 = '11/11/2023 15:56'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-akhdoud-al-hazem-rass/AaT8Kgyg/'
}
$row117 = @{
    6 = 'Al Feiha'
    7 = 0
    8 = 'Al Ettifaq'
    9 = 0
    10 = 2.7
    11 = '04/11/2023 19:13'
    12 = 3.12
    13 = '11/11/2023 15:56'
    14 = 3.21
    15 = '04/11/2023 19:13'
    16 = 3.33
    17 = '11/11/2023 15:56'
    18 = 2.57
    19 = '04/11/2023 19:13'
    20 = 2.36
    21 = '11/11/2023 15:56'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-feiha-al-ettifaq-fc/ddP4LZLn/'
}
Set-RowValues $ws 116 $row117
Set-RowValues $ws 117 $row116

# Swap F:V between row 126 and row 127
$row126 = @{
    6 = 'Al Fateh'
    7 = 0
    8 = 'Al Feiha'
    9 = 1
    10 = 1.97
    11 = '18/11/2023 18:13'
    12 = 1.94
    13 = '25/11/2023 18:34'
    14 = 3.99
    15 = '18/11/2023 18:13'
    16 = 4.11
    17 = '25/11/2023 18:34'
    18 = 3.17
    19 = '18/11/2023 18:13'
    20 = 3.48
    21 = '25/11/2023 18:17'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-fateh-al-feiha/ANB5HrRS/'
}
$row127 = @{
    6 = 'Al Ahli SC'
    7 = 0
    8 = 'Al Shabab'
    9 = 0
    10 = 1.51
    11 = '18/11/2023 19:43'
    12 = 2.06
    13 = '25/11/2023 18:59'
    14 = 4.54
    15 = '18/11/2023 19:43'
    16 = 3.84
    17 = '25/11/2023 18:59'
    18 = 5.05
    19 = '18/11/2023 19:43'
    20 = 3.33
    21 = '25/11/2023 18:59'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ahli-sc-al-shabab/67CzCMlj/'
}
Set-RowValues $ws 126 $row127
Set-RowValues $ws 127 $row126

# Swap F:V between row 135 and row 136
$row135 = @{
    6 = 'Al Raed'
    7 = 2
    8 = 'Al Wehda'
    9 = 0
    10 = 2.39
    11 = '25/11/2023 21:43'
    12 = 3.04
    13 = '02/12/2023 18:58'
    14 = 3.41
    15 = '25/11/2023 21:43'
    16 = 3.52
    17 = '02/12/2023 18:58'
    18 = 2.79
    19 = '25/11/2023 21:43'
    20 = 2.32
    21 = '02/12/2023 18:58'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-raed-al-wehda/jenvWstM/'
}
$row136 = @{
    6 = 'Al Shabab'
    7 = 1
    8 = 'Al Taawon'
    9 = 2
    10 = 2.17
    11 = '25/11/2023 21:43'
    12 = 2.23
    13 = '02/12/2023 18:18'
    14 = 3.41
    15 = '25/11/2023 21:43'
    16 = 3.45
    17 = '02/12/2023 18:59'
    18 = 3.15
    19 = '25/11/2023 21:43'
    20 = 3.26
    21 = '02/12/2023 18:59'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-shabab-al-taawon/b1yEPuJk/'
}
Set-RowValues $ws 135 $row136
Set-RowValues $ws 136 $row135

# Swap F:V between row 140 and row 142
$row140 = @{
    6 = 'Al Taee'
    7 = 1
    8 = 'Al Hilal'
    9 = 2
    10 = 12.91
    11 = '03/12/2023 18:12'
    12 = 20.78
    13 = '08/12/2023 15:59'
    14 = 9.43
    15 = '03/12/2023 18:12'
    16 = 11.18
    17 = '08/12/2023 15:59'
    18 = 1.1
    19 = '03/12/2023 18:12'
    20 = 1.1
    21 = '08/12/2023 15:25'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taee-al-hilal/4dN89pzL/'
}
$row142 = @{
    6 = 'Al Khaleej'
    7 = 3
    8 = 'Abha'
    9 = 1
    10 = 1.87
    11 = '01/12/2023 05:42'
    12 = 1.78
    13 = '08/12/2023 15:34'
    14 = 3.69
    15 = '01/12/2023 05:42'
    16 = 4.22
    17 = '08/12/2023 15:59'
    18 = 3.68
    19 = '01/12/2023 05:42'
    20 = 3.94
    21 = '08/12/2023 15:12'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-khaleej-abha/IamrOqkE/'
}
Set-RowValues $ws 140 $row142
Set-RowValues $ws 142 $row140

# Swap F:V between row 149 and row 150
$row149 = @{
    6 = 'Al Hilal'
    7 = 2
    8 = 'Al Wehda'
    9 = 0
    10 = 1.17
    11 = '09/12/2023 19:13'
    12 = 1.19
    13 = '15/12/2023 15:55'
    14 = 7.44
    15 = '09/12/2023 19:13'
    16 = 7.94
    17 = '15/12/2023 15:55'
    18 = 10.48
    19 = '09/12/2023 19:13'
    20 = 11.76
    21 = '15/12/2023 15:55'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-hilal-al-wehda/Oz1VCJtD/'
}
$row150 = @{
    6 = 'Al Fateh'
    7 = 1
    8 = 'Al Shabab'
    9 = 1
    10 = 2.83
    11 = '08/12/2023 21:43'
    12 = 2.45
    13 = '15/12/2023 15:53'
    14 = 3.49
    15 = '08/12/2023 21:43'
    16 = 3.72
    17 = '15/12/2023 15:57'
    18 = 2.3
    19 = '08/12/2023 21:43'
    20 = 2.72
    21 = '15/12/2023 15:53'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-fateh-al-shabab/SYxIHsSs/'
}
Set-RowValues $ws 149 $row150
Set-RowValues $ws 150 $row149

# Add new rows 154-156, copying the row-153 format first (borders/number formats)
$srcFmt = $ws.Range("A153:V153")
$srcFmt.Copy() | Out-Null
$ws.Range("A154:V154").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A155:V155").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A156:V156").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$newRow154 = @{
    1 = 153
    2 = 'saudi-arabia'
    3 = 'saudi-professional-league'
    4 = '2023-2024'
    5 = 45281.66666666666
    6 = 'Al Khaleej'
    7 = 3
    8 = 'Al Feiha'
    9 = 0
    10 = 2.31
    11 = '16/12/2023 18:13'
    12 = 2
    13 = '21/12/2023 15:59'
    14 = 3.34
    15 = '16/12/2023 18:13'
    16 = 3.53
    17 = '21/12/2023 15:59'
    18 = 2.91
    19 = '16/12/2023 18:13'
    20 = 3.8
    21 = '21/12/2023 15:59'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-khaleej-al-feiha/2RaZBwdJ/'
}
Set-RowValues $ws 154 $newRow154

$newRow155 = @{
    1 = 154
    2 = 'saudi-arabia'
    3 = 'saudi-professional-league'
    4 = '2023-2024'
    5 = 45281.66666666666
    6 = 'Damac'
    7 = 3
    8 = 'Al Taee'
    9 = 0
    10 = 1.81
    11 = '15/12/2023 16:13'
    12 = 1.7
    13 = '21/12/2023 15:51'
    14 = 3.7
    15 = '15/12/2023 16:13'
    16 = 4.05
    17 = '21/12/2023 15:51'
    18 = 3.91
    19 = '15/12/2023 16:13'
    20 = 4.71
    21 = '21/12/2023 15:51'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/damac-al-taee/vBc6EWHH/'
}
Set-RowValues $ws 155 $newRow155

$newRow156 = @{
    1 = 155
    2 = 'saudi-arabia'
    3 = 'saudi-professional-league'
    4 = '2023-2024'
    5 = 45281.79166666666
    6 = 'Al Hilal'
    7 = 7
    8 = 'Abha'
    9 = 0
    10 = 1.05
    11 = '15/12/2023 16:13'
    12 = 1.06
    13 = '21/12/2023 18:55'
    14 = 14.54
    15 = '15/12/2023 16:13'
    16 = 16.25
    17 = '21/12/2023 18:56'
    18 = 16.92
    19 = '15/12/2023 16:13'
    20 = 27.8
    21 = '21/12/2023 18:56'
    22 = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-hilal-abha/AV5kIEIh/'
}
Set-RowValues $ws 156 $newRow156

